# Generate Report for Handoff
# - Flip the localization status from "Handed back: in sync with en-US" to
#   "Ready for handoff" everywhere it is shown (Overview zh-cn/de-de columns
#   and the per-locale Status column), and bump the related handoff
#   timestamps to reflect the newly generated report.
# - The Status column (and the Overview zh-cn/de-de columns that mirror it)
#   shrinks once the shorter text is in place, so re-fit those columns.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

# ---- Overview sheet -------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("G2").Value = "2016-08-29 02:59:04"

# ---- zh-cn sheet ------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("H2").Value = "2016-08-29 02:58:57"

# ---- de-de sheet ------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("H2").Value = "2016-08-29 02:59:04"

# ---- Re-fit the status columns now that the text is shorter -----------
$overview.Columns.Item(5).ColumnWidth = 16.33
$overview.Columns.Item(6).ColumnWidth = 16.33
$zhcn.Columns.Item(3).ColumnWidth = 16.33
$dede.Columns.Item(3).ColumnWidth = 16.33
